$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date values (column B) - whole table shifted back by 4 weeks.
# Leading apostrophe forces these to stay plain text (matches original inlineStr cells).
$dates = @(
    "'2024-12-29",
    "'2025-01-05",
    "'2025-01-12",
    "'2025-01-19",
    "'2025-01-26",
    "'2025-02-02",
    "'2025-02-09",
    "'2025-02-16",
    "'2025-02-23",
    "'2025-03-02",
    "'2025-03-09",
    "'2025-03-16",
    "'2025-03-23",
    "'2025-03-30",
    "'2025-04-06",
    "'2025-04-13"
)

# New MyForecast values (column D)
$forecast = @(12, 14, 10, 11, 11, 12, 12, 13, 12, 12, 12, 11, 12, 12, 10, 12)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $dates[$i]
    $ws1.Cells.Item($row, 4).Value = $forecast[$i]
    # is_holiday_week boolean flag is dropped -> cell becomes blank/empty
    $ws1.Cells.Item($row, 10).Value = ""
}

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'187"
$ws2.Range("B10").Value = "'95"
$ws2.Range("B11").Value = "'47"
$ws2.Range("B12").Value = "'14"
$ws2.Range("B13").Value = "'2025-01-05"
$ws2.Range("B14").Value = "'10"
$ws2.Range("B15").Value = "'2025-04-06"
